$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '71.042.23'
$ws.Range("E2").Value = '  -0.65%  '
Set-TextValue "D3" '3.829.11'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.06%  '
Set-TextValue "D5" '704.09'
$ws.Range("E5").Value = '  +0.48%  '
Set-TextValue "D6" '171.87'
$ws.Range("E6").Value = '  -1.35%  '
Set-TextValue "D7" '3.826.12'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.07%  '
Set-TextValue "D9" '0.525'
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  -0.77%  '
Set-TextValue "D11" '7.41'
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("E13").Value = '  -1.57%  '
Set-TextValue "D14" '36.58'
$ws.Range("E14").Value = '  -0.15%  '
Set-TextValue "D15" '4.475.69'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D16" '71.068.31'
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D17" '3.723.85'
$ws.Range("E17").Value = '  -2.55%  '
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  +0.21%  '
Set-TextValue "D20" '17.37'
$ws.Range("E20").Value = '  -2.43%  '
Set-TextValue "D21" '495.04'
$ws.Range("E21").Value = '  +1.46%  '
Set-TextValue "D22" '10.69'
$ws.Range("E22").Value = '  -4.16%  '
$ws.Range("E23").Value = '  +2.55%  '
Set-TextValue "D24" '85.24'
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("E26").Value = '  +0.62%  '
Set-TextValue "D27" '12.09'
$ws.Range("E27").Value = '  -2.33%  '
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  -1.56%  '
$ws.Range("E31").Value = '  -2.12%  '
$ws.Range("E32").Value = '  -3.55%  '
Set-TextValue "D33" '29.38'
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("E34").Value = '  -3.37%  '
Set-TextValue "D35" '9.19'
$ws.Range("E35").Value = '  -1.48%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue "D36" '3.791.61'
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D37" '1.00'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("E38").Value = '  -1.35%  '
Set-TextValue "D39" '2.33'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("E40").Value = '  +4.14%  '
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("E45").Value = '  +0.43%  '
Set-TextValue "D46" '163.65'
$ws.Range("E46").Value = '  +0.42%  '
Set-TextValue "D47" '428.81'
$ws.Range("E47").Value = '  +4.35%  '
Set-TextValue "D48" '48.87'
$ws.Range("E48").Value = '  +0.24%  '
$ws.Range("E49").Value = '  +0.92%  '
Set-TextValue "D50" '1.37'
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("E51").Value = '  -1.99%  '
